$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(2)

# ---------------------------------------------------------------------------
# 1. Move the "commentID" column (originally column I) so it sits right
#    before "lat" (originally column F). This rotates columns F:I one to
#    the right (old I -> F, old F -> G, old G -> H, old H -> I) for every
#    row in the grid, including the header row.
# ---------------------------------------------------------------------------
$ws.Columns.Item(9).Cut() | Out-Null
$ws.Columns.Item(6).Insert() | Out-Null

# ---------------------------------------------------------------------------
# 2. Insert two new blank rows for the new API calls: "reportComment" goes
#    right before the old row 19 ("deleteMessage2"), and "deleteComment"
#    goes right before the old row 20 ("sendMessage", now at row 21).
# ---------------------------------------------------------------------------
$ws.Rows.Item(19).Insert() | Out-Null
$ws.Rows.Item(21).Insert() | Out-Null

# ---------------------------------------------------------------------------
# 3. Pick up the formatting of a "plain" data row (row 9, "getMessages") for
#    the two new rows, copying columns A:Q so every cell lands on the
#    correct existing style.
# ---------------------------------------------------------------------------
$ws.Range("A9:Q9").Copy() | Out-Null
$ws.Range("A19:Q19").PasteSpecial(-4122) | Out-Null
$ws.Range("A9:Q9").Copy() | Out-Null
$ws.Range("A21:Q21").PasteSpecial(-4122) | Out-Null

# A few cells in the two new rows need the "marked" style (matching the
# darker fill used elsewhere in the grid) instead of the plain one that the
# template row carried.
$ws.Range("F9").Copy() | Out-Null
$ws.Range("H19").PasteSpecial(-4122) | Out-Null
$ws.Range("F9").Copy() | Out-Null
$ws.Range("J19").PasteSpecial(-4122) | Out-Null
$ws.Range("F9").Copy() | Out-Null
$ws.Range("K19").PasteSpecial(-4122) | Out-Null

$ws.Range("F9").Copy() | Out-Null
$ws.Range("H21").PasteSpecial(-4122) | Out-Null
$ws.Range("F9").Copy() | Out-Null
$ws.Range("J21").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 4. Fill in the actual values for the two new rows.
# ---------------------------------------------------------------------------
$ws.Range("B19").Value = "GET"
$ws.Range("C19").Value = "reportComment"
$ws.Range("Q19").Value = "<empty>"

$ws.Range("B21").Value = "GET"
$ws.Range("C21").Value = "deleteComment"

# ---------------------------------------------------------------------------
# 5. Restore the selection shown when the workbook is reopened.
# ---------------------------------------------------------------------------
$ws.Range("O21").Select() | Out-Null
